$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to insert into column B for each data row (2-20), pushing
# existing values one column to the right (B->C, C->D, ... ), dropping
# any value that would overflow past column K.
$newValues = @{
    2  = 0.7916129955631771
    3  = -3.727363316492332
    4  = 0.376932102669816
    5  = 1.207578635508109
    6  = -0.9264868865757077
    7  = 0.3770345820039356
    8  = -0.4275923834192769
    9  = 0.324932645901923
    10 = -0.04071760298358112
    11 = 0.3721869518844864
    12 = -0.1524291232873974
    13 = -1.030518528898312
    14 = 0.4742145784871607
    15 = 0.3556547466179877
    16 = 0.3126006297022321
    17 = 0.3812981176718321
    18 = -0.716162849403934
    19 = 0.506656010950813
    20 = -0.343237405067616
}

$maxCol = 11  # column K

for ($r = 2; $r -le 20; $r++) {
    # Determine the last used column in this row (column B = 2 is the
    # first data column; row 20 starts out with no data at all).
    $lastCol = $ws.Cells.Item($r, $ws.Columns.Count).End(-4159).Column
    if ($lastCol -lt 2) {
        $lastCol = 1
    }

    if ($lastCol -ge 2) {
        # Read the existing values (columns B..lastCol) before overwriting.
        $oldValues = @()
        for ($c = 2; $c -le $lastCol; $c++) {
            $oldValues += , ($ws.Cells.Item($r, $c).Value())
        }

        # Write them back shifted one column to the right, dropping any
        # value that would land past column K (index 11).
        for ($i = 0; $i -lt $oldValues.Count; $i++) {
            $destCol = $i + 3
            if ($destCol -le $maxCol) {
                $ws.Cells.Item($r, $destCol).Value = $oldValues[$i]
            }
        }
    }

    # Insert the new value at the start (column B).
    $ws.Cells.Item($r, 2).Value = $newValues[$r]
}
